$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the textual pin descriptions in column D with plain numeric pin numbers
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 5

# Update the active selection to match the saved workbook state
$ws.Range("D5").Select()
